$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even if it looks numeric or date-like,
# without leaving a lasting NumberFormat/style change on the cell itself.
function Set-TextValue {
    param($rng, [string]$text)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "59.025.80"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "2.500.71"
$ws.Range("E3").Value = "  +1.57%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.25%  "

Set-TextValue $ws.Range("D5") "538.40"
$ws.Range("E5").Value = "  +0.18%  "

Set-TextValue $ws.Range("D6") "143.85"
$ws.Range("E6").Value = "  -3.32%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.65%  "

Set-TextValue $ws.Range("D8") "0.571"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "2.536.01"
$ws.Range("E9").Value = "  +2.77%  "

Set-TextValue $ws.Range("D10") "0.0996"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("E11").Value = "  -0.75%  "

Set-TextValue $ws.Range("D12") "5.50"
$ws.Range("E12").Value = "  +2.95%  "

Set-TextValue $ws.Range("D13") "0.351"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").Value = "2.983.77"
$ws.Range("E14").Value = "  +0.78%  "

Set-TextValue $ws.Range("D15") "23.70"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").Value = "59.074.48"
$ws.Range("E16").Value = "  -2.02%  "

$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "2.525.31"
$ws.Range("E18").Value = "  -2.02%  "

Set-TextValue $ws.Range("D19") "11.26"
$ws.Range("E19").Value = "  +0.13%  "

Set-TextValue $ws.Range("D20") "4.28"
$ws.Range("E20").Value = "  -2.04%  "

Set-TextValue $ws.Range("D21") "323.79"
$ws.Range("E21").Value = "  -0.93%  "

Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  +2.60%  "

Set-TextValue $ws.Range("D23") "5.78"
$ws.Range("E23").Value = "  +0.05%  "

Set-TextValue $ws.Range("D24") "62.08"
$ws.Range("E24").Value = "  +0.38%  "

Set-TextValue $ws.Range("D25") "0.439"
$ws.Range("E25").Value = "  -7.79%  "

$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").Value = "2.627.81"
$ws.Range("E27").Value = "  -1.81%  "

Set-TextValue $ws.Range("D28") "0.991"
$ws.Range("E28").Value = "  -0.32%  "

Set-TextValue $ws.Range("D29") "7.76"
$ws.Range("E29").Value = "  -1.31%  "

Set-TextValue $ws.Range("D30") "6.78"
$ws.Range("E30").Value = "  -3.19%  "

$ws.Range("D31").Value = "0.0₃0774"
$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("E32").Value = "  -3.15%  "

Set-TextValue $ws.Range("D33") "1.20"
$ws.Range("E33").Value = "  -6.92%  "

Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  -0.28%  "

Set-TextValue $ws.Range("D35") "158.08"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("E36").Value = "  +6.10%  "

Set-TextValue $ws.Range("D37") "18.59"
$ws.Range("E37").Value = "  +0.90%  "

Set-TextValue $ws.Range("D38") "4.37"
$ws.Range("E38").Value = "  -5.66%  "

Set-TextValue $ws.Range("D39") "1.61"
$ws.Range("E39").Value = "  -6.22%  "

Set-TextValue $ws.Range("D40") "5.65"
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "36.88"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D42") "301.90"
$ws.Range("E42").Value = "  -4.24%  "

Set-TextValue $ws.Range("D43") "0.819"
$ws.Range("E43").Value = "  -5.32%  "

Set-TextValue $ws.Range("D44") "3.66"
$ws.Range("E44").Value = "  -3.25%  "

Set-TextValue $ws.Range("D45") "0.994"
$ws.Range("E45").Value = "  -0.23%  "

Set-TextValue $ws.Range("D46") "0.602"
$ws.Range("E46").Value = "  +3.72%  "

$ws.Range("E47").Value = "  -0.58%  "

Set-TextValue $ws.Range("D48") "126.17"
$ws.Range("E48").Value = "  +5.12%  "

Set-TextValue $ws.Range("D49") "0.0931"
$ws.Range("E49").Value = "  -1.44%  "

Set-TextValue $ws.Range("D50") "18.70"
$ws.Range("E50").Value = "  -0.75%  "

Set-TextValue $ws.Range("D51") "0.0515"
$ws.Range("E51").Value = "  -1.89%  "
